# A new daily price record was added for "Ajo" (garlic) at Femacal de La
# Calera. In the source data the rows are kept in (roughly) date order, so
# the new record was inserted as row 390, pushing every following row
# (old 390..505) down by one (new 391..506).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 390 - shifts rows 390:505 down to 391:506
$ws.Rows.Item(390).Insert()

# Populate the newly inserted row 390 with the new record's data
$ws.Cells.Item(390, 1).Value  = 3
$ws.Cells.Item(390, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(390, 3).Value  = "Coquimbo"
$ws.Cells.Item(390, 4).Value  = 44784
$ws.Cells.Item(390, 5).Value  = 5
$ws.Cells.Item(390, 6).Value  = 100112003
$ws.Cells.Item(390, 7).Value  = "Ajo"
$ws.Cells.Item(390, 8).Value  = "Chino"
$ws.Cells.Item(390, 9).Value  = "Primera"
$ws.Cells.Item(390, 10).Value = 76
$ws.Cells.Item(390, 11).Value = 24000
$ws.Cells.Item(390, 12).Value = 24500
$ws.Cells.Item(390, 13).Value = 24250
$ws.Cells.Item(390, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(390, 15).Value = "China"
$ws.Cells.Item(390, 16).Value = 2425
$ws.Cells.Item(390, 17).Value = 10
$ws.Cells.Item(390, 18).Value = "Hortaliza"
